$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D; existing D:K data shifts to E:L.
$ws.Columns("D").Insert()

# Copy number formats from the (now-shifted) adjacent column E into the new
# column D for every row that actually carries data, leaving header/spacer
# rows (5, 6, 36, 78, 79) untouched.
$ws.Range("E7:E35,E38:E77,E80:E102").Copy()
$ws.Range("D7:D35,D38:D77,D80:D102").PasteSpecial(-4122)

# Populate the new column D with the latest period's figures.
$ws.Range("D7").Value2 = 43465
$ws.Range("D8").Value2 = 56771000
$ws.Range("D9").Value2 = 37506000
$ws.Range("D10").Value2 = 19265000
$ws.Range("D12").Value2 = "NA"
$ws.Range("D13").Value2 = 0
$ws.Range("D14").Value2 = 0
$ws.Range("D15").Value2 = 3102000
$ws.Range("D17").Value2 = 49531000
$ws.Range("D18").Value2 = 7240000
$ws.Range("D20").Value2 = 0
$ws.Range("D21").Value2 = 10342000
$ws.Range("D22").Value2 = 0
$ws.Range("D23").Value2 = 7240000
$ws.Range("D24").Value2 = -248000
$ws.Range("D25").Value2 = 0
$ws.Range("D26").Value2 = 7488000
$ws.Range("D27").Value2 = 3433000
$ws.Range("D28").Value2 = 0
$ws.Range("D29").Value2 = 0
$ws.Range("D30").Value2 = 0
$ws.Range("D31").Value2 = 0
$ws.Range("D32").Value2 = 0
$ws.Range("D33").Value2 = 3433000
$ws.Range("D34").Value2 = 0
$ws.Range("D35").Value2 = 3433000
$ws.Range("D38").Value2 = 43465
$ws.Range("D41").Value2 = 8390000
$ws.Range("D42").Value2 = 0
$ws.Range("D43").Value2 = 9167000
$ws.Range("D44").Value2 = 6989000
$ws.Range("D45").Value2 = 5508000
$ws.Range("D46").Value2 = 0
$ws.Range("D47").Value2 = 39874000
$ws.Range("D48").Value2 = 151603000
$ws.Range("D49").Value2 = 27577000
$ws.Range("D50").Value2 = 0
$ws.Range("D51").Value2 = 0
$ws.Range("D52").Value2 = 6840000
$ws.Range("D53").Value2 = 0
$ws.Range("D54").Value2 = 256281000
$ws.Range("D57").Value2 = 23740000
$ws.Range("D58").Value2 = 11624000
$ws.Range("D59").Value2 = 0
$ws.Range("D60").Value2 = 0
$ws.Range("D61").Value2 = 108995000
$ws.Range("D62").Value2 = 12236000
$ws.Range("D63").Value2 = 0
$ws.Range("D64").Value2 = 0
$ws.Range("D65").Value2 = 0
$ws.Range("D66").Value2 = 226466000
$ws.Range("D68").Value2 = 0
$ws.Range("D69").Value2 = 0
$ws.Range("D70").Value2 = 4168000
$ws.Range("D71").Value2 = 0
$ws.Range("D72").Value2 = 14551000
$ws.Range("D73").Value2 = 0
$ws.Range("D74").Value2 = 0
$ws.Range("D75").Value2 = 0
$ws.Range("D76").Value2 = 25647000
$ws.Range("D77").Value2 = 0
$ws.Range("D80").Value2 = 43465
$ws.Range("D81").Value2 = 3433000
$ws.Range("D83").Value2 = 3102000
$ws.Range("D84").Value2 = 0
$ws.Range("D85").Value2 = 0
$ws.Range("D86").Value2 = 0
$ws.Range("D87").Value2 = 0
$ws.Range("D88").Value2 = 0
$ws.Range("D89").Value2 = 5158000
$ws.Range("D91").Value2 = -4841000
$ws.Range("D92").Value2 = 0
$ws.Range("D93").Value2 = 0
$ws.Range("D94").Value2 = -19833000
$ws.Range("D96").Value2 = -726000
$ws.Range("D97").Value2 = 0
$ws.Range("D98").Value2 = 0
$ws.Range("D99").Value2 = 0
$ws.Range("D100").Value2 = 18136000
$ws.Range("D101").Value2 = -210000
$ws.Range("D102").Value2 = 3251000
